$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price in column D, Volume(1h) in column E)
$updates = @(
    @{ Row = 2; D = "28.075.40"; DForceText = $false; E = "  +3.38%  " },
    @{ Row = 3; D = "1.688.97"; DForceText = $false; E = "  +0.56%  " },
    @{ Row = 4; D = "0.999"; DForceText = $true; E = "  -0.36%  " },
    @{ Row = 5; D = "216.54"; DForceText = $true; E = "  +1.12%  " },
    @{ Row = 6; D = "0.522"; DForceText = $true; E = "  +0.85%  " },
    @{ Row = 7; D = "0.999"; DForceText = $true; E = "  -0.35%  " },
    @{ Row = 8; D = "24.16"; DForceText = $true; E = "  +7.14%  " },
    @{ Row = 9; D = $null; DForceText = $false; E = "  +2.04%  " },
    @{ Row = 10; D = $null; DForceText = $false; E = "  +0.66%  " },
    @{ Row = 11; D = "0.0886"; DForceText = $true; E = "  -0.56%  " },
    @{ Row = 12; D = $null; DForceText = $false; E = "  +0.45%  " },
    @{ Row = 13; D = "1.686.89"; DForceText = $false; E = "  -0.20%  " },
    @{ Row = 14; D = $null; DForceText = $false; E = "  +0.15%  " },
    @{ Row = 15; D = $null; DForceText = $false; E = "  +0.80%  " },
    @{ Row = 16; D = "66.95"; DForceText = $true; E = "  +0.52%  " },
    @{ Row = 17; D = "250.96"; DForceText = $true; E = "  +6.69%  " },
    @{ Row = 18; D = "28.010.51"; DForceText = $false; E = "  +3.23%  " },
    @{ Row = 19; D = "0.0`u{2083}0743"; DForceText = $false; E = "  +0.72%  " },
    @{ Row = 20; D = "7.67"; DForceText = $true; E = "  -1.93%  " },
    @{ Row = 21; D = "0.999"; DForceText = $true; E = "  -0.31%  " },
    @{ Row = 22; D = $null; DForceText = $false; E = "  +0.09%  " },
    @{ Row = 23; D = $null; DForceText = $false; E = "  +0.65%  " },
    @{ Row = 24; D = $null; DForceText = $false; E = "  -1.43%  " },
    @{ Row = 25; D = "147.40"; DForceText = $true; E = "  +0.42%  " },
    @{ Row = 26; D = "7.37"; DForceText = $true; E = "  -0.25%  " },
    @{ Row = 27; D = "16.48"; DForceText = $true; E = "  +1.16%  " },
    @{ Row = 28; D = $null; DForceText = $false; E = "  +0.51%  " },
    @{ Row = 29; D = $null; DForceText = $false; E = "  -0.32%  " },
    @{ Row = 30; D = $null; DForceText = $false; E = "  +6.83%  " },
    @{ Row = 31; D = "0.0504"; DForceText = $true; E = "  +0.22%  " },
    @{ Row = 32; D = "3.39"; DForceText = $true; E = "  +0.70%  " },
    @{ Row = 33; D = $null; DForceText = $false; E = "  -1.24%  " },
    @{ Row = 34; D = "1.428.70"; DForceText = $false; E = "  -7.19%  " },
    @{ Row = 35; D = $null; DForceText = $false; E = "  -2.16%  " },
    @{ Row = 36; D = "0.946"; DForceText = $true; E = "  +0.56%  " },
    @{ Row = 37; D = $null; DForceText = $false; E = "  -0.25%  " },
    @{ Row = 38; D = $null; DForceText = $false; E = "  -1.90%  " },
    @{ Row = 39; D = "0.0172"; DForceText = $true; E = "  +0.63%  " },
    @{ Row = 40; D = $null; DForceText = $false; E = "  -3.37%  " },
    @{ Row = 41; D = "69.67"; DForceText = $true; E = "  +0.96%  " },
    @{ Row = 42; D = $null; DForceText = $false; E = "  -0.31%  " },
    @{ Row = 43; D = "5.51"; DForceText = $true; E = "  -4.34%  " },
    @{ Row = 44; D = "1.834.77"; DForceText = $false; E = "  +0.35%  " },
    @{ Row = 45; D = "2.23"; DForceText = $true; E = "  -0.75%  " },
    @{ Row = 46; D = "0.797"; DForceText = $true; E = "  +0.85%  " },
    @{ Row = 47; D = $null; DForceText = $false; E = "  +6.02%  " },
    @{ Row = 48; D = "89.43"; DForceText = $true; E = "  -0.51%  " },
    @{ Row = 49; D = "0.0`u{2086}0111"; DForceText = $false; E = "  -0.86%  " },
    @{ Row = 50; D = $null; DForceText = $false; E = "  -0.56%  " },
    @{ Row = 51; D = "7.85"; DForceText = $true; E = "  -4.47%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        if ($u.DForceText) {
            $cellD.NumberFormat = "@"
        }
        $cellD.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
